# Add a new trainee's certificate records (Sayed Marzouk Amin Ali, DSS1401-DSS1408)
# into the 8 empty rows (402-409) that follow the existing data block, mirroring the
# existing layout/formatting used for the previous entrant's 8-course block
# (rows 386-393: same course order/styles, just a new Certificate No + Name).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the previous person's 8-row course block (formatting, course names,
# date, result) straight on top of the blank rows 402:409.
$template = $ws.Range("A386:E393")
$target = $ws.Range("A402:E409")
$template.Copy($target)

# Now overwrite just the Certificate No (col A) and Name (col B) for the new
# trainee. Name is written before Certificate No for each row so the shared
# string table picks up "Sayed Marzouk Amin Ali" ahead of "DSS1401".
$certificateNumbers = @("DSS1401", "DSS1402", "DSS1403", "DSS1404", "DSS1405", "DSS1406", "DSS1407", "DSS1408")
$traineeName = "Sayed Marzouk Amin Ali"

for ($i = 0; $i -lt $certificateNumbers.Length; $i++) {
    $row = 402 + $i
    $ws.Cells.Item($row, 2).Value = $traineeName
    $ws.Cells.Item($row, 1).Value = $certificateNumbers[$i]
}

# Match the author's final on-screen view: scrolled down to row 393 and the
# active cell sitting on D411.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 393
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D411").Select()
